# Remove the two "Medium" sub-bullets that are no longer relevant:
#   - "Where to show error message generated when user tries to delete a
#      currency and fails"
#   - "Check referential integrity when we delete a record and modify the
#      behavior if necessary, that is decide prevent deleting a record if
#      it has child records or delete child records as well"
#
# Both paragraphs immediately follow the "Medium" heading paragraph and
# immediately precede the "Remove all id fields from attr_accessible"
# paragraph, so locate them by their text and delete the whole range they
# occupy (including their paragraph marks).

$d = $word.ActiveDocument

$startText = "Where to show error message generated when user tries to delete a currency and fails"
$endText   = "Check referential integrity when we delete a record and modify the behavior if necessary, that is decide prevent deleting a record if it has child records or delete child records as well"

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -match [regex]::Escape($startText)) {
        $startPara = $p
    }
    if ($t -match [regex]::Escape($endText)) {
        $endPara = $p
    }
}

$range = $d.Range($startPara.Range.Start, $endPara.Range.End)
$range.Delete()
